$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "RateCompare": fix LV3 label typos (drop stray underscore before v#)
# ---------------------------------------------------------------------
$rc = $wb.Worksheets.Item("RateCompare")
$rc.Range("A3").Value = "LV3_130v1"
$rc.Range("A4").Value = "LV3_130v2"
$rc.Range("A5").Value = "LV3_200v1"
$rc.Range("A6").Value = "LV3_200v2"

# ---------------------------------------------------------------------
# Insert a new data row (LV11_glc20) above the Wierckx rows, pushing the
# two Wierckx rows down by one.
# ---------------------------------------------------------------------
$rc.Rows.Item(7).Insert()

$rc.Range("A7").Value = "LV11_glc20"
$rc.Range("B7").Value = "6-35"
$rc.Range("C7").Value = "lin"
$rc.Range("D7").Value = 0.17
$rc.Range("E7").Value = 0.02
$rc.Range("F7").Value = 0.55
$rc.Range("G7").Value = 0.12
$rc.Range("H7").Value = "EX_glc__D_e"

# Normalise formatting across the affected block (rows 7-9) so every cell
# shares the same (default) cell style.
$rc.Range("A7:H9").Font.Name = "Arial"
$rc.Range("A7:H9").Font.Size = 10

$rc.Application.GoTo($rc.Range("H7"))

# ---------------------------------------------------------------------
# Sheet "Metadata": small cosmetic touch-up + selection move
# ---------------------------------------------------------------------
$md = $wb.Worksheets.Item("Metadata")
$md.Range("E1").Font.Name = "Arial"
$md.Range("E1").Font.Size = 10

$md.Application.Union($md.Range("E2"), $rc.Range("H7")).Select()
